$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (Excel row number) and column letter
$updates = @{
    2  = @{ B = -0.1079215458242564; C = 0.4990845831178726; D = 0.6310303664915219; E = 0.7943741980273037; F = 0.7948402182622575; G = 51 }
    3  = @{ B = 0.3713778311963193;  C = 0.7365167446170454; D = 1.455560046373499;  E = 1.206465932537467;  F = 1.159538328638186;  G = 50 }
    4  = @{ B = 0.4609105965926396;  C = 1.00885168430549;   D = 3.578481368703621;  E = 1.891687439484552;  F = 1.85369055727329;   G = 49 }
    5  = @{ B = 0.4265731976396117;  C = 1.138403783358868;  D = 4.557394820086173;  E = 2.134805569621312;  F = 2.113888451448179;  G = 48 }
    6  = @{ B = 0.2577980780148568;  C = 0.9475793177038206; D = 3.943416922645604;  E = 1.985803847978346;  F = 1.990286097306333;  G = 47 }
    7  = @{ B = 0.2936268802263522;  C = 0.999724168007887;  D = 4.905488059981674;  E = 2.214833641604189;  F = 2.227333515476834;  G = 35 }
    8  = @{ B = 0.2339287874485977;  C = 1.007627164027908;  D = 5.003647654014131;  E = 2.236883469028758;  F = 2.258072718037386;  G = 34 }
    9  = @{ B = 0.1807689603703594;  C = 1.497342948823048;  D = 9.511094063079014;  E = 3.084006171050735;  F = 3.173455165340056;  G = 17 }
    10 = @{ B = -0.6878977726065058; C = 1.16036325979617;   D = 7.183952125666254;  E = 2.680289560041275;  F = 2.730638667684675 }
    11 = @{ B = 0.1218082912510491;  C = 0.4641168854429287; D = 0.2644604851195395; E = 0.5142572168861994; F = 0.5585955886216332 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
